$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "purchased"
$ws.Range("C1").Select() | Out-Null
